# Apply updated cryptocurrency Price (D) and Volume(1h) (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.722.82"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.859.88"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").Value = "'1.020"
$ws.Range("E4").Value = "  -0.99%  "
$ws.Range("D5").Value = "'320.85"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'1.018"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").Value = "'0.4371"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").Value = "'0.3777"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("D9").Value = "'0.07410"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "'0.8842"
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("D11").Value = "'21.56"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").Value = "1.858.95"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "'6.758"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "'5.495"
$ws.Range("E14").Value = "  -0.72%  "
$ws.Range("D15").Value = "'0.07142"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "'87.84"
$ws.Range("D17").Value = "'1.023"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "'0.000009026"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "'15.46"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "27.720.46"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "'5.286"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'11.16"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("D24").Value = "2.088.40"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'2.037"
$ws.Range("E25").Value = "  +6.04%  "
$ws.Range("D26").Value = "'157.14"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "'18.69"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "'5.426"
$ws.Range("E28").Value = "  +2.32%  "
$ws.Range("D29").Value = "'1.989"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").Value = "'121.34"
$ws.Range("E30").Value = "  +3.49%  "
$ws.Range("D31").Value = "'0.09047"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "'1.217"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "'0.7701"
$ws.Range("E33").Value = "  +0.75%  "
$ws.Range("D34").Value = "'3.036"
$ws.Range("E34").Value = "  +5.27%  "
$ws.Range("D35").Value = "'4.560"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'1.020"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "'1.138"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").Value = "'0.01978"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "'0.05309"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "'2.877"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("D41").Value = "'0.5181"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "'6.959"
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("D43").Value = "'0.1678"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "'8.709"
$ws.Range("E44").Value = "  +2.09%  "
$ws.Range("D45").Value = "'10.79"
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("D46").Value = "'110.33"
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").Value = "'1.711"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").Value = "'1.019"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").Value = "'0.06476"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "'1.845"
$ws.Range("E51").Value = "  -0.85%  "
